$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.156683206558228
$ws.Range("B1").Value = 2.383946180343628
$ws.Range("D1").Value = 2.389873504638672
$ws.Range("E1").Value = 1.223288774490356
